$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The song title previously had the artist baked into the title string
# ("<Title> - John Mayer"). Split it: trim the title back down to just the
# song name, and record "John Mayer" as the Authors value for that row
# (it was previously "Unknown").
$oldTitle = "You_re No One Til Someone Lets You Down - John Mayer"
$newTitle = "You_re No One Til Someone Lets You Down"
$author   = "John Mayer"

$cell = $ws.Cells.Find($oldTitle)
$row = $cell.Row

$ws.Cells.Item($row, 2).Value = $newTitle
$ws.Cells.Item($row, 3).Value = $author

# Column A's stored width tracked the longest title in column B, which was
# this row's (it had " - John Mayer" - 13 characters - tacked onto the end).
# Now that the title has been trimmed, narrow column A by that same 13
# characters: ~53.7 -> ~40.7.
# (This host quantizes ColumnWidth to the nearest 1/6th of a character on
# save; 39.8333... is the input that lands on the closest achievable stored
# width, 40.6667, to the intended 40.7109375.)
$ws.Columns.Item(1).ColumnWidth = 39.833333333333336
